$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.503597122302158
$ws.Range("C2").Value = 0.563106796116505
$ws.Range("D2").Value = 0.573643410852713
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 0.491499227202473

$ws.Range("B3").Value = 0.855345911949686
$ws.Range("C3").Value = 0.735135135135135
$ws.Range("D3").Value = 0.824427480916031
$ws.Range("E3").Value = 0.734939759036145
$ws.Range("F3").Value = 0.537890044576523

$ws.Range("B4").Value = 0.771428571428571
$ws.Range("C4").Value = 0.752747252747253
$ws.Range("D4").Value = 0.783783783783784
$ws.Range("E4").Value = 0.738853503184713
$ws.Range("F4").Value = 0.550375939849624

$ws.Range("B5").Value = 0.884353741496599
$ws.Range("C5").Value = 0.809782608695652
$ws.Range("D5").Value = 0.859259259259259
$ws.Range("E5").Value = 0.865030674846626
$ws.Range("F5").Value = 0.609422492401216
